$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 390 ("「賢くバカであれ」" entry), shifting all subsequent rows up by one.
$ws.Rows.Item(390).Delete()
